$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-24 Wednesday", "2025-09-25 Thursday"),
    @("793×2=1586", "387×2=774"),
    @("922×5=4610", "736×3=2208"),
    @("977×7=6839", "682×9=6138"),
    @("585×5=2925", "805×6=4830"),
    @("199×5=995", "999×9=8991"),
    @("969×3=2907", "191×4=764"),
    @("359×9=3231", "760×6=4560"),
    @("501×8=4008", "117×8=936"),
    @("635×6=3810", "149×4=596"),
    @("458×7=3206", "870×9=7830"),
    @("580×8=4640", "889×4=3556"),
    @("336×9=3024", "358×5=1790"),
    @("893×4=3572", "487×6=2922"),
    @("557×4=2228", "247×4=988"),
    @("551×7=3857", "431×7=3017"),
    @("936×6=5616", "568×9=5112"),
    @("406×5=2030", "561×9=5049"),
    @("207×4=828", "141×8=1128"),
    @("979×6=5874", "706×2=1412"),
    @("406×4=1624", "858×8=6864"),
    @("904×3=2712", "231×3=693"),
    @("690×4=2760", "306×4=1224"),
    @("416×7=2912", "384×5=1920"),
    @("334×2=668", "760×8=6080"),
    @("149×7=1043", "794×3=2382")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
